$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.617.23"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "2.614.15"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "511.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.85%  "
$ws.Range("D9").Value = "2.627.26"
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.60%  "
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.346"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").Value = "3.073.56"
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").Value = "60.541.36"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "2.626.66"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "352.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").Value = "0.0₃0844"
$ws.Range("E28").Value = "  -2.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "151.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("E35").Value = "  -1.57%  "
$ws.Range("E36").Value = "  -1.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.891"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.847"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.16%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "293.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.625"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.100"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("E47").Value = "  -4.42%  "
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").Value = "2.005.88"
$ws.Range("E51").Value = "  -2.54%  "
